$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B: new genome sequence (rows 1-20)
$bValues = @(0,12,4,14,8,18,13,16,15,2,9,11,6,5,19,10,17,1,3,7)
for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}

# Column D: fitness / penalty improvement values
$ws.Range("D1").Value = 109.2448836041059
$ws.Range("D2").Value = 74.52689319035321

# Row 21, column B: last generation fit
$ws.Range("B21").Value = 0.8015459867609394
